# feat: add 2022-Q3 data
#
# 1. Insert a new "2022-Q3" sheet right after "总计", shifting all the
#    later quarter sheets down by one tab position.
# 2. Populate the new sheet with the 2022-Q3 holdings detail (one fund row).
# 3. Insert a new row on "总计" for the 2022-Q3 summary (count=1, value=0.01)
#    and keep the running index column sequential.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)        # "总计"

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $totalSheet)
$q3.Name = "2022-Q3"

# NOTE: worksheet variables resolve by tab position, so grab the template
# sheet reference AFTER the insert above (it has just been pushed from
# position 2 to position 3) - fetching it earlier would silently re-point
# at the freshly inserted "2022-Q3" sheet instead.
$templateSheet = $wb.Worksheets.Item("2022-Q2")   # detail-sheet layout template

# Pull the header-row (B1:H1) formatting + the index-column (A2) formatting
# from the existing template sheet so the new sheet matches the others.
$templateSheet.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$templateSheet.Range("A2").Copy()
$q3.Range("A2").PasteSpecial(-4122)

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0

# Fund code / size / position columns are stored as text in this workbook
# (leading zeros in the code, fixed-format decimals) - force text storage.
$q3.Range("B2").NumberFormat = "@"
$q3.Range("D2:G2").NumberFormat = "@"

$q3.Range("B2").Value = "013242"
$q3.Range("C2").Value = "北信瑞丰优势行业股票"
$q3.Range("D2").Value = "0.69"
$q3.Range("E2").Value = "92.79"
$q3.Range("F2").Value = "1.55"
$q3.Range("G2").Value = "0.0107"
$q3.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2) Add the matching summary row on "总计"
# ---------------------------------------------------------------------
$totalSheet.Rows.Item(2).Insert()

# Clear formatting the insert may have inherited from the header row so the
# new data row looks like a normal data row (no bold / no border).
$totalSheet.Range("B3").Copy()
$totalSheet.Range("B2").PasteSpecial(-4122)
$totalSheet.Range("C3:D3").Copy()
$totalSheet.Range("C2:D2").PasteSpecial(-4122)

$totalSheet.Range("B2").Value = "2022-Q3"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

# Re-number the running index column (A) sequentially, 0-based, including
# the newly appended row at the bottom.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5
$totalSheet.Range("A8").Value = 6

Write-Output "2022-Q3 sheet + summary row added"
